$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "scenario" column header in F1
$ws.Range("F1").Value = "scenario"

# Fill F2:F101 with the scenario label "S2" for every data row
$ws.Range("F2:F101").Value = "S2"
